$wb = $excel.ActiveWorkbook

# The sheet "Killed get_players progressbar" is copied to the end of the
# workbook (same mechanism Excel uses for "Move or Copy... Create a copy"),
# producing a brand-new sheet that keeps the original layout/formulas/values.
$src = $wb.Worksheets.Item("Killed get_players progressbar")
[void]$src.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Removed Index from Player "

# The player index columns were re-measured without the old index column,
# so the raw timing entries on the new sheet get updated values (the
# dependent AVERAGE/SUM/diff formulas recompute automatically).
$newSheet.Range("B2").Value = 8.1549999999999994
$newSheet.Range("C2").Value = 8.407
$newSheet.Range("D2").Value = 8.2449999999999992

$newSheet.Range("B3").Value = 7.415
$newSheet.Range("C3").Value = 7.6440000000000001
$newSheet.Range("D3").Value = 7.5049999999999999

$newSheet.Range("B4").Value = 0.73699999999999999
$newSheet.Range("C4").Value = 0.746
$newSheet.Range("D4").Value = 0.73599999999999999

# Keep the new sheet active/selected (tab highlighted, D5 selected like the
# original "Killed get_players progressbar" sheet used to be), and reset the
# selection on the original sheet back to its full used range.
$src.PageSetup.Orientation = 1
[void]$src.Range("A1:G6").Select()
[void]$newSheet.Activate()
[void]$newSheet.Range("D5").Select()
